# "Created TC2 - InvalidLogin"
#
# Renames the empty "Sheet3" test-case sheet to "InvalidLogin", fills it in
# with the same username/password header row used by the "ValidLogin" sheet
# plus a sample invalid credential pair ("abc" / "xyz"), and switches the
# workbook's active tab / selection state from ValidLogin over to the new
# InvalidLogin sheet.

$wb = $excel.ActiveWorkbook

$wsValid   = $wb.Worksheets.Item(2)   # "ValidLogin"
$wsInvalid = $wb.Worksheets.Item(3)   # currently "Sheet3"

# Rename Sheet3 -> InvalidLogin (TC2)
$wsInvalid.Name = "InvalidLogin"

# Header row (bold), matching the style used on the ValidLogin sheet
$wsInvalid.Range("A1").Value = "username"
$wsInvalid.Range("B1").Value = "password"
$wsInvalid.Range("A1:B1").Font.Bold = $true

# Sample invalid login data
$wsInvalid.Range("A2").Value = "abc"
$wsInvalid.Range("B2").Value = "xyz"

# ValidLogin is no longer the selected/active tab; its selection becomes A1:B2
$wsValid.Activate()
$wsValid.Range("A1:B2").Select()

# InvalidLogin becomes the active tab, with D5 selected
$wsInvalid.Activate()
$wsInvalid.Range("D5").Select()
